$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new timelog entry as row 28 (mirrors the formatting of row 27:
# style 3 / numFmt-16 on column A, style 4 / wrapped-Arial on column B).
$ws.Range("A28").Value = "3/12, 4 hrs"
$ws.Range("B28").Value = "Making things pretty, documenting, organizing, adding final touches"

# Copy the formatting (styles/borders) from the row above so the new row
# matches the rest of the table (s="3" on A28, s="4" on B28).
$ws.Range("A27:B27").Copy()
$ws.Range("A28:B28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the row height used by the other wrapped-text rows of this size.
$ws.Rows.Item(28).RowHeight = 41.4

# Update the view: the sheet is scrolled down a couple more rows and the
# active/selected cell moves on to the next (empty) row of the log.
$ws.Range("B29").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
